$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 258, pushing existing rows 258-331 down to 264-337
$ws.Rows.Item(258).Resize(6).Insert()

# Constant columns shared by every data row in this sheet
$A = 8
$B = "Terminal La Palmera de La Serena"
$C = "Coquimbo"
$E = 4
$F = "Fruta"
$G = 100103
$H = "Frutos de hueso (carozo)"
$I = 100103004
$J = "Durazno"

# New rows data: D, K, L, M, N, O, P, Q, R, S, T
$rows = @(
    @{ R=258; D=44588; K="Andross";      L="Especial"; M=20; N=360000; O=365000; P=362500; Q="`$/bins (400 kilos)"; Rg="Región de O'Higgins"; S=906; T=400 },
    @{ R=259; D=44588; K="Andross";      L="Primera";  M=20; N=320000; O=325000; P=322500; Q="`$/bins (400 kilos)"; Rg="Región de O'Higgins"; S=806; T=400 },
    @{ R=260; D=44588; K="Andross";      L="Segunda";  M=16; N=260000; O=265000; P=262500; Q="`$/bins (400 kilos)"; Rg="Región de O'Higgins"; S=656; T=400 },
    @{ R=261; D=44588; K="Elegant Lady"; L="Especial"; M=16; N=370000; O=375000; P=372500; Q="`$/bins (400 kilos)"; Rg="Región de O'Higgins"; S=931; T=400 },
    @{ R=262; D=44588; K="Elegant Lady"; L="Primera";  M=20; N=330000; O=335000; P=332500; Q="`$/bins (400 kilos)"; Rg="Región de O'Higgins"; S=831; T=400 },
    @{ R=263; D=44588; K="Elegant Lady"; L="Segunda";  M=20; N=300000; O=305000; P=302500; Q="`$/bins (400 kilos)"; Rg="Región de O'Higgins"; S=756; T=400 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $A
    $ws.Cells.Item($r, 2).Value = $B
    $ws.Cells.Item($r, 3).Value = $C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $E
    $ws.Cells.Item($r, 6).Value = $F
    $ws.Cells.Item($r, 7).Value = $G
    $ws.Cells.Item($r, 8).Value = $H
    $ws.Cells.Item($r, 9).Value = $I
    $ws.Cells.Item($r, 10).Value = $J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.Rg
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
}

Write-Host "Done"
